# Fruta / hortaliza, semanal
# Insert two new weekly rows of data right above the current row 30,
# pushing all subsequent rows (old 30..105) down to (32..107).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 30 (shifts existing rows 30-105 down to 32-107)
$ws.Rows("30:31").Insert()

# --- New row 30 ---
$ws.Range("A30").Value = 4
$ws.Range("B30").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C30").Value = "Los Lagos"
$ws.Range("D30").Value = 44565
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103002
$ws.Range("J30").Value = "Ciruela"
$ws.Range("K30").Value = "Black Amber"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 600
$ws.Range("N30").Value = 18000
$ws.Range("O30").Value = 18500
$ws.Range("P30").Value = 18250
$ws.Range("Q30").Value = "`$/caja 15 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 1217
$ws.Range("T30").Value = 15

# --- New row 31 ---
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44565
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = "Frutos de hueso (carozo)"
$ws.Range("I31").Value = 100103002
$ws.Range("J31").Value = "Ciruela"
$ws.Range("K31").Value = "Black Amber"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 250
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("Q31").Value = "`$/caja 15 kilos granel"
$ws.Range("R31").Value = "Región de O'Higgins"
$ws.Range("S31").Value = 1067
$ws.Range("T31").Value = 15
